$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - update short URL entry
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = "https://www.baeldung.com/properties-with-spring"
$ws.Range("C1").Value = "dnc29h4a"

# Row 2 - update short URL entry
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "www.google.com"
$ws.Range("C2").Value = "7nm7eid9"

# Remove the old extra rows (3 and 4) that are no longer present
$ws.Rows("3:4").Delete()

# Resize the Long URL / Short code columns
$ws.Columns("B").ColumnWidth = 24.666666666666668
$ws.Columns("C").ColumnWidth = 37.666666666666664

# Update the active selection to C1
[void]$ws.Range("C1").Select()
